# Auto-generated update of cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'64.930.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.09%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.515.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.09%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'587.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.63%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'133.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.92%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'3.514.52"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.09%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  +0.03%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  -0.71%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = "'  +2.07%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'7.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.36%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.384"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.51%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'4.110.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.14%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'27.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.83%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'  +0.22%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("E16").Value = "'  +0.81%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'3.513.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.29%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'64.962.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.04%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'10.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.18%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'14.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.16%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'5.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.34%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'390.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.73%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.574"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.08%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'74.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.49%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'3.657.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.19%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  -0.08%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "'  -2.41%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  +8.25%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'7.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.33%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  -0.14%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'2.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.92%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'8.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.92%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'3.523.30"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.82%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'24.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.73%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'  +0.05%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  +2.08%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'5.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +4.50%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  +2.27%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'169.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.92%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'6.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.96%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.0804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.54%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  -0.38%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'26.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.30%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'42.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.67%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'  +4.16%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "'  +0.04%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'4.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.21%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'  +1.24%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'2.495.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.94%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = "'  -0.24%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("B51").Value = "'SuiNetwork"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.894"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.50%  "
$ws.Range("E51").Style = "Normal"

